# Edit script: 
# 1) Replace the text "Job: Delivery Man" with "Button response".
# 2) Move the "_GoBack" bookmark from the "Should accept data" paragraph
#    to the (now renamed) "Button response" paragraph.

$d = $word.ActiveDocument

# --- Step 1: remove the existing _GoBack bookmark (if present) ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- Step 2: change the run text "Job: Delivery Man" -> "Button response" ---
$d.Content.Find.Execute("Job: Delivery Man", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Button response", 2) | Out-Null

# --- Step 3: re-add the _GoBack bookmark right after the new text ---
$range = $d.Content
$range.Find.Execute("Button response", $true, $false, $false, $false, $false,
                     $true, 1, $false, "", 0) | Out-Null
$range.Collapse(0)
$d.Bookmarks.Add("_GoBack", $range)
